$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: fix the "Th" + bookmark + "is is a sub-project..." split so the
# sentence reads naturally, and split "Date & Time" into its own Heading3
# paragraph (moving the "_GoBack" bookmark to sit before "Time").
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(15)   # "This is a sub-project of the super-project Circle Language Spec."
$p2 = $d.Paragraphs(16)   # "Time" (Heading3)
$combined = $d.Range($p1.Range.Start, $p2.Range.End)

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# rsid attributes below mirror the originating paragraphs/runs (000F2F0D /
# 009D1238 for the sub-project sentence, 00FF5F4B for the Time heading, and
# the untouched italic run's 00B33380) so only the genuinely-edited content
# picks up fresh (rsid-less) runs, same as the rest of this document's
# plain <w:r> runs.
$newBody = @"
<w:p w:rsidR="000F2F0D" w:rsidRDefault="000F2F0D" w:rsidP="009D1238" $xmlNs>
  <w:pPr>
    <w:ind w:left="426"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">This is a sub-project of the super-project </w:t>
  </w:r>
  <w:r w:rsidR="00B33380" w:rsidRPr="00B33380">
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>Circle Language Spec</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p w:rsidR="00FF5F4B" w:rsidRDefault="00FF5F4B" w:rsidP="00FF5F4B" $xmlNs>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Date &amp; </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t>Time</w:t>
  </w:r>
</w:p>
"@

$combined.InsertXML($newBody)

# ---------------------------------------------------------------------
# Change 2 & 3: reorder the w:attr children of the two date smart tags'
# smartTagPr from Month/Day/Year to Year/Day/Month.
# ---------------------------------------------------------------------
$dateRangePara = $d.Paragraphs(17)  # holds both "March 5, 2008" and "April 13, 2008" smart tags

# Same rsid attributes as the original paragraph/runs (00FF5F4B / 00916C2E /
# 009D1238) - only the w:attr order inside each smartTagPr changes.
$dateBody = '<w:p w:rsidR="00FF5F4B" w:rsidRDefault="00916C2E" w:rsidP="009D1238" ' + $xmlNs + '><w:pPr><w:ind w:left="426"/></w:pPr>' + `
  '<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">' + `
    '<w:smartTagPr>' + `
      '<w:attr w:name="Year" w:val="2008"/>' + `
      '<w:attr w:name="Day" w:val="5"/>' + `
      '<w:attr w:name="Month" w:val="3"/>' + `
    '</w:smartTagPr>' + `
    '<w:r><w:t xml:space="preserve">March 5, </w:t></w:r>' + `
    '<w:r w:rsidR="00FF5F4B"><w:t>2008</w:t></w:r>' + `
  '</w:smartTag>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidR="00FF5F4B"><w:t xml:space="preserve">&#8211; </w:t></w:r>' + `
  '<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">' + `
    '<w:smartTagPr>' + `
      '<w:attr w:name="Year" w:val="2008"/>' + `
      '<w:attr w:name="Day" w:val="13"/>' + `
      '<w:attr w:name="Month" w:val="4"/>' + `
    '</w:smartTagPr>' + `
    '<w:r><w:t>April 13, 2008</w:t></w:r>' + `
  '</w:smartTag>' + `
'</w:p>'

$dateRangePara.Range.InsertXML($dateBody)

Write-Output "edit applied"
